# Auto-generated update of market price columns (H-N) across all Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")


# --- ALC sheet ---
# Row 19
$ws_ALC.Range("H19").Value = 417.45
$ws_ALC.Range("I19").Value = 298.4
$ws_ALC.Range("J19").Value = 536.5
$ws_ALC.Range("K19").Value = 298.4
$ws_ALC.Range("L19").Value = 536.5
$ws_ALC.Range("M19").Value = -123.4
$ws_ALC.Range("N19").Value = -886.5
# Row 41
$ws_ALC.Range("H41").Value = 1233.9333
$ws_ALC.Range("I41").Value = 339.25
$ws_ALC.Range("J41").Value = 2256.4285
$ws_ALC.Range("K41").Value = 339.25
$ws_ALC.Range("L41").Value = 2256.4285
$ws_ALC.Range("M41").Value = 100.75
$ws_ALC.Range("N41").Value = -3136.4285
# Row 98
$ws_ALC.Range("H98").Value = 2839.639
$ws_ALC.Range("I98").Value = 2885.4849
$ws_ALC.Range("K98").Value = 2885.4849
$ws_ALC.Range("M98").Value = -1387.4849
# Row 113
$ws_ALC.Range("H113").Value = 2299
$ws_ALC.Range("J113").Value = 2749
$ws_ALC.Range("L113").Value = 2749
$ws_ALC.Range("N113").Value = -9257
# Row 122
$ws_ALC.Range("H122").Value = 2839.639
$ws_ALC.Range("I122").Value = 2885.4849
$ws_ALC.Range("K122").Value = 8656.4547
$ws_ALC.Range("M122").Value = -6206.4547
# Row 132
$ws_ALC.Range("H132").Value = 9670.823
$ws_ALC.Range("I132").Value = 9670.823
$ws_ALC.Range("K132").Value = 29012.469
$ws_ALC.Range("M132").Value = -26482.469
# Row 138
$ws_ALC.Range("H138").Value = 1479061.8
$ws_ALC.Range("I138").Value = 13123.75
$ws_ALC.Range("J138").Value = 1930119.5
$ws_ALC.Range("K138").Value = 39371.25
$ws_ALC.Range("L138").Value = 5790358.5
$ws_ALC.Range("M138").Value = -34231.25
$ws_ALC.Range("N138").Value = -5800638.5

# --- ARM sheet ---
# Row 74
$ws_ARM.Range("H74").Value = 197819.17
$ws_ARM.Range("I74").Value = 271821.72
$ws_ARM.Range("J74").Value = 3562.5
$ws_ARM.Range("K74").Value = 271821.72
$ws_ARM.Range("L74").Value = 3562.5
$ws_ARM.Range("M74").Value = -270947.72
$ws_ARM.Range("N74").Value = -5310.5
# Row 77
$ws_ARM.Range("H77").Value = 197819.17
$ws_ARM.Range("I77").Value = 271821.72
$ws_ARM.Range("J77").Value = 3562.5
$ws_ARM.Range("K77").Value = 1359108.6
$ws_ARM.Range("L77").Value = 17812.5
$ws_ARM.Range("M77").Value = -1354740.6
$ws_ARM.Range("N77").Value = -26548.5
# Row 110
$ws_ARM.Range("H110").Value = 1939.6389
$ws_ARM.Range("I110").Value = 1058.6333
$ws_ARM.Range("K110").Value = 1058.6333
$ws_ARM.Range("M110").Value = 986.3667
# Row 122
$ws_ARM.Range("H122").Value = 4071.3704
$ws_ARM.Range("I122").Value = 3842.8076
$ws_ARM.Range("J122").Value = 10014
$ws_ARM.Range("K122").Value = 11528.4228
$ws_ARM.Range("L122").Value = 30042
$ws_ARM.Range("M122").Value = -9078.4228
$ws_ARM.Range("N122").Value = -34942
# Row 132
$ws_ARM.Range("H132").Value = 3472.8518
$ws_ARM.Range("J132").Value = 4265.5
$ws_ARM.Range("L132").Value = 12796.5
$ws_ARM.Range("N132").Value = -17856.5

# --- BSM sheet ---
# Row 94
$ws_BSM.Range("H94").Value = 55556380
$ws_BSM.Range("I94").Value = 62500772
$ws_BSM.Range("J94").Value = 1234
$ws_BSM.Range("K94").Value = 62500772
$ws_BSM.Range("L94").Value = 1234
$ws_BSM.Range("M94").Value = -62500321
$ws_BSM.Range("N94").Value = -2136
# Row 132
$ws_BSM.Range("H132").Value = 96327.664
$ws_BSM.Range("J132").Value = 96327.664
$ws_BSM.Range("L132").Value = 96327.664
$ws_BSM.Range("N132").Value = -106447.664

# --- CRP sheet ---
# Row 16
$ws_CRP.Range("H16").Value = 1523.85
$ws_CRP.Range("I16").Value = 1321.4706
$ws_CRP.Range("J16").Value = 2670.6667
$ws_CRP.Range("K16").Value = 1321.4706
$ws_CRP.Range("L16").Value = 2670.6667
$ws_CRP.Range("M16").Value = -1034.4706
$ws_CRP.Range("N16").Value = -3244.6667
# Row 31
$ws_CRP.Range("H31").Value = 5258.3853
$ws_CRP.Range("I31").Value = 4323.295
$ws_CRP.Range("J31").Value = 6888.1143
$ws_CRP.Range("K31").Value = 4323.295
$ws_CRP.Range("L31").Value = 6888.1143
$ws_CRP.Range("M31").Value = -4028.295
$ws_CRP.Range("N31").Value = -7478.1143
# Row 34
$ws_CRP.Range("H34").Value = 5258.3853
$ws_CRP.Range("I34").Value = 4323.295
$ws_CRP.Range("J34").Value = 6888.1143
$ws_CRP.Range("K34").Value = 4323.295
$ws_CRP.Range("L34").Value = 6888.1143
$ws_CRP.Range("M34").Value = -4121.295
$ws_CRP.Range("N34").Value = -7292.1143
# Row 58
$ws_CRP.Range("H58").Value = 3039.6875
$ws_CRP.Range("I58").Value = 1762.7778
$ws_CRP.Range("K58").Value = 1762.7778
$ws_CRP.Range("M58").Value = -1559.7778
# Row 92
$ws_CRP.Range("H92").Value = 61663.332
$ws_CRP.Range("J92").Value = 61663.332
$ws_CRP.Range("L92").Value = 61663.332
$ws_CRP.Range("N92").Value = -66655.33199999999
# Row 96
$ws_CRP.Range("H96").Value = 9000
$ws_CRP.Range("J96").Value = 9000
$ws_CRP.Range("L96").Value = 9000
$ws_CRP.Range("N96").Value = -14492
# Row 107
$ws_CRP.Range("H107").Value = 1689
$ws_CRP.Range("I107").Value = 1318
$ws_CRP.Range("K107").Value = 1318
$ws_CRP.Range("M107").Value = 602
# Row 113
$ws_CRP.Range("H113").Value = 1523.85
$ws_CRP.Range("I113").Value = 1321.4706
$ws_CRP.Range("J113").Value = 2670.6667
$ws_CRP.Range("K113").Value = 1321.4706
$ws_CRP.Range("L113").Value = 2670.6667
$ws_CRP.Range("M113").Value = 848.5293999999999
$ws_CRP.Range("N113").Value = -7010.6667
# Row 122
$ws_CRP.Range("H122").Value = 2795.925
$ws_CRP.Range("I122").Value = 2800.2666
$ws_CRP.Range("J122").Value = 2782.9
$ws_CRP.Range("K122").Value = 8400.799800000001
$ws_CRP.Range("L122").Value = 8348.700000000001
$ws_CRP.Range("M122").Value = -5950.799800000001
$ws_CRP.Range("N122").Value = -13248.7
# Row 132
$ws_CRP.Range("H132").Value = 11115135
$ws_CRP.Range("I132").Value = 12823855
$ws_CRP.Range("J132").Value = 8455
$ws_CRP.Range("K132").Value = 38471565
$ws_CRP.Range("L132").Value = 25365
$ws_CRP.Range("M132").Value = -38469035
$ws_CRP.Range("N132").Value = -30425
# Row 134
$ws_CRP.Range("H134").Value = 6030.467
$ws_CRP.Range("I134").Value = 5911.241
$ws_CRP.Range("J134").Value = 6246.5625
$ws_CRP.Range("K134").Value = 17733.723
$ws_CRP.Range("L134").Value = 18739.6875
$ws_CRP.Range("M134").Value = -15198.723
$ws_CRP.Range("N134").Value = -23809.6875
# Row 136
$ws_CRP.Range("H136").Value = 3039.6875
$ws_CRP.Range("I136").Value = 1762.7778
$ws_CRP.Range("K136").Value = 5288.3334
$ws_CRP.Range("M136").Value = -2738.3334

# --- CUL sheet ---
# Row 14
$ws_CUL.Range("H14").Value = 558.63635
$ws_CUL.Range("I14").Value = 558.63635
$ws_CUL.Range("K14").Value = 1675.90905
$ws_CUL.Range("M14").Value = -1502.90905
# Row 56
$ws_CUL.Range("H56").Value = 7233.75
$ws_CUL.Range("I56").Value = 7233.75
$ws_CUL.Range("K56").Value = 7233.75
$ws_CUL.Range("M56").Value = -6703.75
# Row 133
$ws_CUL.Range("H133").Value = 15997.5
$ws_CUL.Range("I133").Value = 12996.667
$ws_CUL.Range("K133").Value = 38990.001
$ws_CUL.Range("M133").Value = -33930.001
# Row 134
$ws_CUL.Range("H134").Value = 1693.16
$ws_CUL.Range("I134").Value = 1201.619
$ws_CUL.Range("K134").Value = 3604.857
$ws_CUL.Range("M134").Value = 1465.143

# --- GSM sheet ---
# Row 5
$ws_GSM.Range("H5").Value = 10000
$ws_GSM.Range("J5").Value = 0
$ws_GSM.Range("L5").Value = 0
$ws_GSM.Range("N5").ClearContents()
# Row 102
$ws_GSM.Range("H102").Value = 2957.111
$ws_GSM.Range("I102").Value = 2949.3076
$ws_GSM.Range("K102").Value = 2949.3076
$ws_GSM.Range("M102").Value = -1327.3076

# --- LTW sheet ---
# Row 40
$ws_LTW.Range("H40").Value = 8055.343
$ws_LTW.Range("J40").Value = 10000
$ws_LTW.Range("L40").Value = 10000
$ws_LTW.Range("N40").Value = -10272
# Row 55
$ws_LTW.Range("H55").Value = 533.36365
$ws_LTW.Range("I55").Value = 535.2222
$ws_LTW.Range("J55").Value = 525
$ws_LTW.Range("K55").Value = 535.2222
$ws_LTW.Range("L55").Value = 525
$ws_LTW.Range("M55").Value = -362.2222
$ws_LTW.Range("N55").Value = -871
# Row 61
$ws_LTW.Range("H61").Value = 2615.8696
$ws_LTW.Range("I61").Value = 1568.7333
$ws_LTW.Range("K61").Value = 1568.7333
$ws_LTW.Range("M61").Value = -1366.7333
# Row 113
$ws_LTW.Range("H113").Value = 2615.8696
$ws_LTW.Range("I113").Value = 1568.7333
$ws_LTW.Range("K113").Value = 1568.7333
$ws_LTW.Range("M113").Value = 601.2666999999999
# Row 122
$ws_LTW.Range("H122").Value = 5198.8184
$ws_LTW.Range("I122").Value = 5439.4443
$ws_LTW.Range("J122").Value = 4116
$ws_LTW.Range("K122").Value = 16318.3329
$ws_LTW.Range("L122").Value = 12348
$ws_LTW.Range("M122").Value = -13868.3329
$ws_LTW.Range("N122").Value = -17248
# Row 136
$ws_LTW.Range("H136").Value = 6390.2085
$ws_LTW.Range("I136").Value = 4880.227
$ws_LTW.Range("K136").Value = 14640.681
$ws_LTW.Range("M136").Value = -12090.681
# Row 140
$ws_LTW.Range("H140").Value = 78189.75
$ws_LTW.Range("J140").Value = 78189.75
$ws_LTW.Range("L140").Value = 78189.75
$ws_LTW.Range("N140").Value = -88549.75
# Row 141
$ws_LTW.Range("H141").Value = 112663.836
$ws_LTW.Range("J141").Value = 112663.836
$ws_LTW.Range("L141").Value = 112663.836
$ws_LTW.Range("N141").Value = -123023.836

# --- WVR sheet ---
# Row 122
$ws_WVR.Range("H122").Value = 17861710
$ws_WVR.Range("I122").Value = 3993.0908
$ws_WVR.Range("K122").Value = 11979.2724
$ws_WVR.Range("M122").Value = -9529.2724
# Row 126
$ws_WVR.Range("H126").Value = 2680.5264
$ws_WVR.Range("I126").Value = 2507.8125
$ws_WVR.Range("K126").Value = 7523.4375
$ws_WVR.Range("M126").Value = -5053.4375
# Row 132
$ws_WVR.Range("H132").Value = 4633266
$ws_WVR.Range("J132").Value = 3612.7856
$ws_WVR.Range("L132").Value = 10838.3568
$ws_WVR.Range("N132").Value = -15898.3568
# Row 135
$ws_WVR.Range("H135").Value = 57525.43
$ws_WVR.Range("J135").Value = 57525.43
$ws_WVR.Range("L135").Value = 57525.43
$ws_WVR.Range("N135").Value = -67665.42999999999
